$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "Through 2022-03-08"

# 2. Update the header text for column B (cutoff date changed from March 07 to March 08)
$ws.Range("B1").Value = "March 2022 (through March 08)"

# 3. Insert a new row before row 55 (shifts "Clearing"..."Wrigleyville" down by one row)
$ws.Rows.Item(55).Insert()

# 4. Populate the newly inserted row 55 with the "Boystown" data
$ws.Range("A55").Value = "Boystown"
$ws.Range("A55").Font.Bold = $true
$ws.Range("A55").HorizontalAlignment = -4108
$ws.Range("A55").VerticalAlignment = -4160
$ws.Range("A55").Borders.LineStyle = 1
$ws.Range("B55").Value = 1

# 5. Apply isolated single-cell updates across the rest of the sheet
$ws.Range("K4").Value = 1
$ws.Range("N5").Value = 4
$ws.Range("N12").Value = 1
$ws.Range("W12").Value = 1
$ws.Range("K14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("W15").Value = 2
$ws.Range("H17").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("Q41").Value = 1
